$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Ativacao date change (keep as literal text, not auto-converted to a date serial)
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2021"
$ws.Range("B9").Copy()
$ws.Range("B8").PasteSpecial(-4122)

$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "01/01/2021"
$ws.Range("C9").Copy()
$ws.Range("C8").PasteSpecial(-4122)

# 2. Add English "Objectives" translation into existing row 11 (A11 already has the label)
$ws.Range("B13").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("C13").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("B11").Value = 'Introduce the fundamental concepts of management science and organization selttings.'
$ws.Range("C11").Value = 'Introduce the fundamental concepts of management science and organization selttings.'

# 3. Docentes responsaveis - new teacher
$ws.Range("B13").Value = '11079086 - Herlandí de Souza Andrade'
$ws.Range("C13").Value = '11079086 - Herlandí de Souza Andrade'

# 4. Programa resumido - merge lines (remove inner newline)
$ws.Range("B14").Value = '1. Áreas de Atuação da Administração.2. Estrutura organizacional.'
$ws.Range("C14").Value = '1. Áreas de Atuação da Administração.2. Estrutura organizacional.'

# 5. Add English "Short syllabus" translation into existing row 15
$ws.Range("B14").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("B15").Value = '1. Management Practice Areas. 2. Organizational structure'
$ws.Range("C15").Value = '1. Management Practice Areas. 2. Organizational structure'

# 6. Programa - merge lines (remove inner newline)
$ws.Range("B16").Value = '1. Noções básicas de Marketing, Finanças e Recursos Humanos. 2. Diferentes configurações de organização.'
$ws.Range("C16").Value = '1. Noções básicas de Marketing, Finanças e Recursos Humanos. 2. Diferentes configurações de organização.'

# 7. Add English "Syllabus" translation into existing row 17
$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("C16").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("B17").Value = '1. Basic notions of Marketing, Finance and Human Resources.2. Different organization settings.'
$ws.Range("C17").Value = '1. Basic notions of Marketing, Finance and Human Resources.2. Different organization settings.'

# 8. Metodo
$ws.Range("B19").Value = 'Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras.'
$ws.Range("C19").Value = 'Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras.'

# 9. Criterio
$ws.Range("B20").Value = 'Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas'
$ws.Range("C20").Value = 'Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas'

# 10. Norma de recuperacao
$ws.Range("B21").Value = 'NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação'
$ws.Range("C21").Value = 'NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação'

# 11. Bibliografia
$ws.Range("B22").Value = 'Chiavenato, I. Gestão de Pessoas. 4 ed. São Paulo: Manole, 2014.Chiavenato, I. Recursos Humanos: o capital humano das organizações. 10 ed. Rio de Janeiro, Campus, 2015.ROBBINS, S. P.; DECENZO, D. A.; WOLTER, R. Fundamentos de Gestão de Pessoas. São Paulo, saraiva, 2013.KOTLER, P. - ARMSTRONG, G. Princípios De Marketing. 15 ed. São Paulo: Pearson, 2014.KOTLER, P.; KELLER, K. L. Administração De Marketing. 15 ed. São Paulo: Pearson, 2019.CHIAVENATO, I. Introdução À Teoria Geral da Administração. 9 ed. São Paulo: Manole, 2014. MAXIMIANO, A. C. A. Teoria Geral da Administração: da Revolução Urbana À Revolução Digital. 8 ed. São Paulo: Atlas, 2017.GUERRINI, F. M.; ESCRIÇÃO FILHO, E.; ROSIM, D. Administração Para Engenheiros. Rio de Janeiro: Campus, 2016.CHIAVENATO, I. Administração Para Não Administradores: a Gestão de Negócios Ao Alcance de Todos. 2 ed. São Paulo: Manole, 2011.SILVA, M. M. L. Administração para Estudantes e Profissionais de Áreas Técnicas. São Paulo: Brasport, 2018.GITMAN, L. J. - ZUTTER, C. J. Princípios de Administração Financeira. 14 ed. São Paulo: Perason, 2017.GROPPELLI, A. A.; NIKBAKHT, E. Administração Financeira. 3 ed. São Paulo: Saraiva, 2010.MARCOUSÉ, I.; SURRIDGE, M.; GILLESPIE, A. Finanças. São Paulo: Saraiva, 2013.BOLMAN, L.G.; DEAL, T.E. Reframing organizations. San Francisco, John Wiley, 2013KOTLER, P.. O Marketing sem segredos. 1 ed. Porto Alegre. Bookman, 2005MINTZBERG, H. Criando organizações eficazes. 2 ed. São Paulo, Atlas, 2006.MORGAN, G. Imagens da organização. São Paulo, Atlas, 1996.'
$ws.Range("C22").Value = 'Chiavenato, I. Gestão de Pessoas. 4 ed. São Paulo: Manole, 2014.Chiavenato, I. Recursos Humanos: o capital humano das organizações. 10 ed. Rio de Janeiro, Campus, 2015.ROBBINS, S. P.; DECENZO, D. A.; WOLTER, R. Fundamentos de Gestão de Pessoas. São Paulo, saraiva, 2013.KOTLER, P. - ARMSTRONG, G. Princípios De Marketing. 15 ed. São Paulo: Pearson, 2014.KOTLER, P.; KELLER, K. L. Administração De Marketing. 15 ed. São Paulo: Pearson, 2019.CHIAVENATO, I. Introdução À Teoria Geral da Administração. 9 ed. São Paulo: Manole, 2014. MAXIMIANO, A. C. A. Teoria Geral da Administração: da Revolução Urbana À Revolução Digital. 8 ed. São Paulo: Atlas, 2017.GUERRINI, F. M.; ESCRIÇÃO FILHO, E.; ROSIM, D. Administração Para Engenheiros. Rio de Janeiro: Campus, 2016.CHIAVENATO, I. Administração Para Não Administradores: a Gestão de Negócios Ao Alcance de Todos. 2 ed. São Paulo: Manole, 2011.SILVA, M. M. L. Administração para Estudantes e Profissionais de Áreas Técnicas. São Paulo: Brasport, 2018.GITMAN, L. J. - ZUTTER, C. J. Princípios de Administração Financeira. 14 ed. São Paulo: Perason, 2017.GROPPELLI, A. A.; NIKBAKHT, E. Administração Financeira. 3 ed. São Paulo: Saraiva, 2010.MARCOUSÉ, I.; SURRIDGE, M.; GILLESPIE, A. Finanças. São Paulo: Saraiva, 2013.BOLMAN, L.G.; DEAL, T.E. Reframing organizations. San Francisco, John Wiley, 2013KOTLER, P.. O Marketing sem segredos. 1 ed. Porto Alegre. Bookman, 2005MINTZBERG, H. Criando organizações eficazes. 2 ed. São Paulo, Atlas, 2006.MORGAN, G. Imagens da organização. São Paulo, Atlas, 1996.'

# 12. Remove Requisitos rows (23 and 24)
$ws.Rows("23:24").Delete()
